$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.523361921310425
$ws.Range("B1").Value = 3.339061975479126
$ws.Range("C1").Value = 2.975239515304565
$ws.Range("D1").Value = 2.508692264556885
$ws.Range("E1").Value = 1.650665760040283
